$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing table:
# Serie / Cupo / Monto demandado / Monto total adjudicado /
# Monto adjudicado bancos y soc. financieras / Monto adjudicado AFP y otras / Tasa
$dates = @("12-10-2021", "13-10-2021", "14-10-2021", "15-10-2021")

$rows = @(
    @{ B = 2000000; C = 2607000; D = 2000000; E = 1525000; F = 475000;  G = 2.5  },
    @{ B = 2000000; C = 3435000; D = 2000000; E = 1680000; F = 320000;  G = 2.5  },
    @{ B = 3500000; C = 5255000; D = 4550000; E = 3805000; F = 745000;  G = 2.75 },
    @{ B = 4000000; C = 7028000; D = 6000000; E = 4502000; F = 1498000; G = 2.75 }
)

$startRow = 196

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $dateCell = $ws.Cells.Item($r, 1)
    $dateText = $dates[$i]

    if ($dateText -eq "12-10-2021") {
        # "12-10-2021" looks like a valid month-day-year date (Dec 10, 2021)
        # to Excel's auto-recognition, so it would otherwise be silently
        # converted into a date serial number. Force it to be entered as
        # plain text, then restore the default "Normal" style so the cell
        # ends up looking exactly like the other plain-text date cells in
        # the column (the other three dates, e.g. "13-10-2021", are not
        # valid dates and are therefore already kept as text automatically).
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $dateText
        $dateCell.Style = "Normal"
    } else {
        $dateCell.Value = $dateText
    }

    $ws.Cells.Item($r, 2).Value = $rows[$i].B
    $ws.Cells.Item($r, 3).Value = $rows[$i].C
    $ws.Cells.Item($r, 4).Value = $rows[$i].D
    $ws.Cells.Item($r, 5).Value = $rows[$i].E
    $ws.Cells.Item($r, 6).Value = $rows[$i].F
    $ws.Cells.Item($r, 7).Value = $rows[$i].G
}
